# Update the NATMI TPM-derived metrics (columns E..T) on rows 2-11 of the
# active sheet with newly recomputed values from the updated TPM pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @('E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T')

$data = @{}
$data[2]  = @(1,0.3333333333333333,0.155562,0.466686,0.04499457894025669,0.05066525794134161,2,0.6666666666666666,0.4472413333333334,1.341724,0.4361236687690723,0.4361236687690722,0.06957375629600002,0.6261638066640001,0.01962320084214439,0.02209631817250927)
$data[3]  = @(1,0.3333333333333333,0.155562,0.466686,0.04499457894025669,0.05066525794134161,2,0.6666666666666666,0.5782506666666666,1.734752,0.5638763312309278,0.5638763312309278,0.08995383020799999,0.8095844718719999,0.02537137809811231,0.02856893976883234)
$data[4]  = @(3,1,1.748562333333333,5.245687,0.5057522141597954,0.5694923030357509,2,0.6666666666666666,0.4472413333333334,1.341724,0.4361236687690723,0.4361236687690722,0.7820293493764445,7.038264144388001,0.2205705111274515,0.2483690725357)
$data[5]  = @(3,1,1.748562333333333,5.245687,0.5057522141597954,0.5694923030357509,2,0.6666666666666666,0.5782506666666666,1.734752,0.5638763312309278,0.5638763312309278,1.011107334958222,9.099966014624,0.285181703032344,0.321123230500051)
$data[6]  = @(2,0.6666666666666666,0.2140146666666667,0.6420440000000001,0.06190136288879069,0.06970280846155817,2,0.6666666666666666,0.4472413333333334,1.341724,0.4361236687690723,0.4361236687690722,0.09571620487288891,0.8614458438560002,0.02699664948486509,0.03039904454976268)
$data[7]  = @(2,0.6666666666666666,0.2140146666666667,0.6420440000000001,0.06190136288879069,0.06970280846155817,2,0.6666666666666666,0.5782506666666666,1.734752,0.5638763312309278,0.5638763312309278,0.1237541236764444,1.113787113088,0.03490471340392561,0.0393037639117955)
$data[8]  = @(2,1,1.1608855,2.321771,0.3357732239901092,0.2520605430540592,2,0.6666666666666666,0.4472413333333334,1.341724,0.4361236687690723,0.4361236687690722,0.5191959788673334,3.115175873204,0.1464386503209859,0.109929568788661)
$data[9]  = @(2,1,1.1608855,2.321771,0.3357732239901092,0.2520605430540592,2,0.6666666666666666,0.5782506666666666,1.734752,0.5638763312309278,0.5638763312309278,0.6712828142986665,4.027696885791999,0.1893345736691233,0.1421309742653982)
$data[10] = @(1,0.3333333333333333,0.1783253333333333,0.534976,0.05157862002104791,0.05807908750729007,2,0.6666666666666666,0.4472413333333334,1.341724,0.4361236687690723,0.4361236687690722,0.07975445984711112,0.7177901386240001,0.02249465699362534,0.02532966472243933)
$data[11] = @(1,0.3333333333333333,0.1783253333333333,0.534976,0.05157862002104791,0.05807908750729007,2,0.6666666666666666,0.5782506666666666,1.734752,0.5638763312309278,0.5638763312309278,0.1031167428835555,0.9280506859519999,0.02908396302742258,0.03274942278485073)

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "$($cols[$i])$row"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
